$wb = $excel.ActiveWorkbook

$wsTools = $wb.Worksheets.Item("Tools")
$wsSources = $wb.Worksheets.Item("Sources")

# Update the "Sources" table: the link in B11 now points to the same URL but
# with updated link text ("Department for Education (gov.uk)" instead of
# "Employer representative bodies"). The A11 description text is unchanged.
$wsSources.Range("B11").Value = "<a href='https://www.gov.uk/government/publications/designated-employer-representative-bodies/notice-of-designated-employer-representative-bodies'>Department for Education (gov.uk)</a>"

# Update the last-selected cell on each sheet to match the saved view state.
$wsTools.Activate() | Out-Null
$wsTools.Range("C10").Select() | Out-Null

$wsSources.Activate() | Out-Null
$wsSources.Range("F8").Select() | Out-Null
